# Auto-generated script applying scheduled-runner value updates
# to the Famfrit_Profits workbook (profit calc sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4661.8887
$ws.Range("I111").Value = 4699.087
$ws.Range("J111").Value = 4448
$ws.Range("K111").Value = 14097.261
$ws.Range("L111").Value = 13344
$ws.Range("M111").Value = -11030.261
$ws.Range("N111").Value = -19478

$ws.Range("H113").Value = 4487.9414
$ws.Range("I113").Value = 3390.5454
$ws.Range("J113").Value = 6499.8335
$ws.Range("K113").Value = 3390.5454
$ws.Range("L113").Value = 6499.8335
$ws.Range("M113").Value = -136.5454
$ws.Range("N113").Value = -13007.8335

$ws.Range("H135").Value = 16668659
$ws.Range("I135").Value = 1544.3334
$ws.Range("K135").Value = 13899.0006
$ws.Range("M135").Value = -11364.0006

$ws.Range("H141").Value = 3252
$ws.Range("I141").Value = 3252
$ws.Range("K141").Value = 9756
$ws.Range("M141").Value = -4576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1877.7646
$ws.Range("I2").Value = 1869.3125
$ws.Range("K2").Value = 1869.3125
$ws.Range("M2").Value = -1756.3125

$ws.Range("H61").Value = 29414632
$ws.Range("I61").Value = 41668012
$ws.Range("K61").Value = 41668012
$ws.Range("M61").Value = -41667800

$ws.Range("H110").Value = 10785.216
$ws.Range("I110").Value = 13217.896
$ws.Range("J110").Value = 1966.75
$ws.Range("K110").Value = 13217.896
$ws.Range("L110").Value = 1966.75
$ws.Range("M110").Value = -11172.896
$ws.Range("N110").Value = -6056.75

$ws.Range("H116").Value = 1877.7646
$ws.Range("I116").Value = 1869.3125
$ws.Range("K116").Value = 1869.3125
$ws.Range("M116").Value = 424.6875

$ws.Range("H122").Value = 4981.25
$ws.Range("J122").Value = 4981.25
$ws.Range("L122").Value = 14943.75
$ws.Range("N122").Value = -19843.75

$ws.Range("H132").Value = 38464916
$ws.Range("I132").Value = 2944.3333
$ws.Range("K132").Value = 8832.999899999999
$ws.Range("M132").Value = -6302.999899999999

$ws.Range("H136").Value = 29414632
$ws.Range("I136").Value = 41668012
$ws.Range("K136").Value = 125004036
$ws.Range("M136").Value = -125001486

$ws.Range("H138").Value = 119999.5
$ws.Range("J138").Value = 119999.5
$ws.Range("L138").Value = 119999.5
$ws.Range("N138").Value = -130279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1877.7646
$ws.Range("I3").Value = 1869.3125
$ws.Range("K3").Value = 1869.3125
$ws.Range("M3").Value = -1755.3125

$ws.Range("H80").Value = 2078.2632
$ws.Range("J80").Value = 4145.625
$ws.Range("L80").Value = 4145.625
$ws.Range("N80").Value = -6141.625

$ws.Range("H83").Value = 2078.2632
$ws.Range("J83").Value = 4145.625
$ws.Range("L83").Value = 20728.125
$ws.Range("N83").Value = -30712.125

$ws.Range("H107").Value = 3200.3513
$ws.Range("I107").Value = 2112.1924
$ws.Range("K107").Value = 2112.1924
$ws.Range("M107").Value = -192.1923999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 952.44446
$ws.Range("I16").Value = 968.3333
$ws.Range("K16").Value = 968.3333
$ws.Range("M16").Value = -681.3333

$ws.Range("H31").Value = 20412896
$ws.Range("I31").Value = 3909.0386
$ws.Range("J31").Value = 43483924
$ws.Range("K31").Value = 3909.0386
$ws.Range("L31").Value = 43483924
$ws.Range("M31").Value = -3614.0386
$ws.Range("N31").Value = -43484514

$ws.Range("H34").Value = 20412896
$ws.Range("I34").Value = 3909.0386
$ws.Range("J34").Value = 43483924
$ws.Range("K34").Value = 3909.0386
$ws.Range("L34").Value = 43483924
$ws.Range("M34").Value = -3707.0386
$ws.Range("N34").Value = -43484328

$ws.Range("H58").Value = 1870.2354
$ws.Range("J58").Value = 2317.818
$ws.Range("L58").Value = 2317.818
$ws.Range("N58").Value = -2723.818

$ws.Range("H80").Value = 46998.332
$ws.Range("J80").Value = 46998.332
$ws.Range("L80").Value = 46998.332
$ws.Range("N80").Value = -49244.332

$ws.Range("H83").Value = 46998.332
$ws.Range("J83").Value = 46998.332
$ws.Range("L83").Value = 140994.996
$ws.Range("N83").Value = -152226.996

$ws.Range("H97").Value = 66196.25
$ws.Range("J97").Value = 66196.25
$ws.Range("L97").Value = 66196.25
$ws.Range("N97").Value = -68178.25

$ws.Range("H105").Value = 13030.667
$ws.Range("I105").Value = 1628.6
$ws.Range("K105").Value = 1628.6
$ws.Range("M105").Value = 118.4000000000001

$ws.Range("H113").Value = 952.44446
$ws.Range("I113").Value = 968.3333
$ws.Range("K113").Value = 968.3333
$ws.Range("M113").Value = 1201.6667

$ws.Range("H136").Value = 1870.2354
$ws.Range("J136").Value = 2317.818
$ws.Range("L136").Value = 6953.454000000001
$ws.Range("N136").Value = -12053.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1110.1111
$ws.Range("I114").Value = 497.75
$ws.Range("J114").Value = 1600
$ws.Range("K114").Value = 1493.25
$ws.Range("L114").Value = 4800
$ws.Range("M114").Value = 1760.75
$ws.Range("N114").Value = -11308

$ws.Range("H132").Value = 2225038.2
$ws.Range("I132").Value = 2175.6667
$ws.Range("J132").Value = 3177693.8
$ws.Range("K132").Value = 19581.0003
$ws.Range("L132").Value = 28599244.2
$ws.Range("M132").Value = -17051.0003
$ws.Range("N132").Value = -28604304.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 21094
$ws.Range("I20").Value = 15000
$ws.Range("J20").Value = 22617.5
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 22617.5
$ws.Range("M20").Value = -14755
$ws.Range("N20").Value = -23107.5

$ws.Range("H113").Value = 2511.8
$ws.Range("I113").Value = 1673.4762
$ws.Range("K113").Value = 1673.4762
$ws.Range("M113").Value = 496.5237999999999

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 716.3333
$ws.Range("I16").Value = 716.3333
$ws.Range("K16").Value = 716.3333
$ws.Range("M16").Value = -546.3333

$ws.Range("H25").Value = 27535
$ws.Range("I25").Value = 7505
$ws.Range("K25").Value = 7505
$ws.Range("M25").Value = -7275

$ws.Range("H35").Value = 719.8
$ws.Range("I35").Value = 719.8
$ws.Range("K35").Value = 719.8
$ws.Range("M35").Value = -383.8

$ws.Range("H55").Value = 692.8421
$ws.Range("J55").Value = 1115.375
$ws.Range("L55").Value = 1115.375
$ws.Range("N55").Value = -1461.375

$ws.Range("H61").Value = 4700.875
$ws.Range("I61").Value = 4241.3335
$ws.Range("J61").Value = 6079.5
$ws.Range("K61").Value = 4241.3335
$ws.Range("L61").Value = 6079.5
$ws.Range("M61").Value = -4039.3335
$ws.Range("N61").Value = -6483.5

$ws.Range("H68").Value = 2409.9565
$ws.Range("J68").Value = 2975
$ws.Range("L68").Value = 2975
$ws.Range("N68").Value = -4473

$ws.Range("H71").Value = 2409.9565
$ws.Range("J71").Value = 2975
$ws.Range("L71").Value = 14875
$ws.Range("N71").Value = -22363

$ws.Range("H113").Value = 4700.875
$ws.Range("I113").Value = 4241.3335
$ws.Range("J113").Value = 6079.5
$ws.Range("K113").Value = 4241.3335
$ws.Range("L113").Value = 6079.5
$ws.Range("M113").Value = -2071.3335
$ws.Range("N113").Value = -10419.5

$ws.Range("H134").Value = 80864
$ws.Range("J134").Value = 80864
$ws.Range("L134").Value = 80864
$ws.Range("N134").Value = -91004

$ws.Range("H136").Value = 1839.7593
$ws.Range("I136").Value = 1873
$ws.Range("K136").Value = 5619
$ws.Range("M136").Value = -3069

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 29000
$ws.Range("J29").Value = 50000
$ws.Range("L29").Value = 50000
$ws.Range("N29").Value = -50580

$ws.Range("H32").Value = 176333.33
$ws.Range("I32").Value = 176333.33
$ws.Range("K32").Value = 176333.33
$ws.Range("M32").Value = -176016.33

$ws.Range("H37").Value = 37247.5
$ws.Range("I37").Value = 45000
$ws.Range("J37").Value = 29495
$ws.Range("K37").Value = 45000
$ws.Range("L37").Value = 29495
$ws.Range("M37").Value = -44797
$ws.Range("N37").Value = -29901

$ws.Range("H40").Value = 14582.5
$ws.Range("I40").Value = 7666.6665
$ws.Range("K40").Value = 7666.6665
$ws.Range("M40").Value = -7517.6665

$ws.Range("H43").Value = 19999
$ws.Range("J43").Value = 19999
$ws.Range("L43").Value = 19999
$ws.Range("N43").Value = -20297

$ws.Range("H107").Value = 485.88235
$ws.Range("I107").Value = 397.36365
$ws.Range("K107").Value = 1192.09095
$ws.Range("M107").Value = 727.90905

$ws.Range("H132").Value = 4299.2856
$ws.Range("I132").Value = 4267.8184
$ws.Range("K132").Value = 12803.4552
$ws.Range("M132").Value = -10273.4552
